# Add the two new use-case values for Alyssa Robinson's row (row 3)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "RemoveMessage"
$ws.Range("B3").Value = "DeleteMessage"

# Update the selected/active cell to B3 to match the saved view state
$ws.Range("B3").Select()
